$p = $ppt.ActivePresentation

# --- Slide 2: "Hybrid data model" -> "Graph+TimeSeries Hybrid data model" ---
# (both occurrences, in the "Content Placeholder 3" shape)
$s2 = $p.Slides.Item(2)
$sh2 = $s2.Shapes.Item(4)
$tr2 = $sh2.TextFrame.TextRange

$para1 = $tr2.Paragraphs(1, 1)
$para1.Text = "X"
$para1.Text = "Graph+TimeSeries Hybrid data model"

$para2 = $tr2.Paragraphs(2, 1)
$para2.Text = "X"
$para2.Text = "Graph+TimeSeries Hybrid data model"

# --- Slide 3: remove the stray "Wor" paragraph ---
$s3 = $p.Slides.Item(3)
$sh3 = $s3.Shapes.Item(1)
$tr3 = $sh3.TextFrame.TextRange

$wor = $tr3.Paragraphs(4, 1)
$wor.Delete()
